$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 135
$ws.Range("I2").Value = 304
$ws.Range("J2").Value = 1234
$ws.Range("K2").Value = 5
$ws.Range("L2").Value = 317
$ws.Range("M2").Value = 14
$ws.Range("N2").Value = 223
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 3
$ws.Range("Q2").Value = 2
$ws.Range("R2").Value = 14
$ws.Range("T2").Value = 215
$ws.Range("U2").Value = 12
$ws.Range("V2").Value = 1941
$ws.Range("W2").Value = 2
$ws.Range("X2").Value = 1910
$ws.Range("Y2").Value = 4
$ws.Range("Z2").Value = 26
$ws.Range("AA2").Value = 9
